$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap / rotate existing match rows (F:V) to reflect corrected fixtures ---
# Row 61
$ws.Range("F61").Value = "Hillerod"
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = "Horsens"
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 2.44
$ws.Range("K61").Value = "24/09/2023 16:13"
$ws.Range("L61").Value = 2.98
$ws.Range("M61").Value = "29/09/2023 18:32"
$ws.Range("N61").Value = 3.52
$ws.Range("O61").Value = "24/09/2023 16:13"
$ws.Range("P61").Value = 3.55
$ws.Range("Q61").Value = "29/09/2023 18:32"
$ws.Range("R61").Value = 2.66
$ws.Range("S61").Value = "24/09/2023 16:13"
$ws.Range("T61").Value = 2.34
$ws.Range("U61").Value = "29/09/2023 18:32"
$ws.Range("V61").Value = "https://www.betexplorer.com/football/denmark/1st-division/hillerod-horsens/jqvSYQd2/"

# Row 62
$ws.Range("F62").Value = "Vendsyssel"
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = "Aalborg"
$ws.Range("I62").Value = 3
$ws.Range("J62").Value = 3.87
$ws.Range("K62").Value = "23/09/2023 17:13"
$ws.Range("L62").Value = 4.85
$ws.Range("M62").Value = "29/09/2023 18:46"
$ws.Range("N62").Value = 3.84
$ws.Range("O62").Value = "23/09/2023 17:13"
$ws.Range("P62").Value = 4.09
$ws.Range("Q62").Value = "29/09/2023 18:47"
$ws.Range("R62").Value = 1.85
$ws.Range("S62").Value = "23/09/2023 17:13"
$ws.Range("T62").Value = 1.67
$ws.Range("U62").Value = "29/09/2023 18:45"
$ws.Range("V62").Value = "https://www.betexplorer.com/football/denmark/1st-division/vendsyssel-ff-aalborg/0E7GucVE/"

# Row 68
$ws.Range("F68").Value = "Horsens"
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = "Aalborg"
$ws.Range("I68").Value = 4
$ws.Range("J68").Value = 3.66
$ws.Range("K68").Value = "29/09/2023 18:13"
$ws.Range("L68").Value = 3.88
$ws.Range("M68").Value = "06/10/2023 18:53"
$ws.Range("N68").Value = 3.86
$ws.Range("O68").Value = "29/09/2023 18:13"
$ws.Range("P68").Value = 3.72
$ws.Range("Q68").Value = "06/10/2023 18:54"
$ws.Range("R68").Value = 1.84
$ws.Range("S68").Value = "29/09/2023 18:13"
$ws.Range("T68").Value = 1.93
$ws.Range("U68").Value = "06/10/2023 18:54"
$ws.Range("V68").Value = "https://www.betexplorer.com/football/denmark/1st-division/horsens-aalborg/SUv8VC77/"

# Row 69
$ws.Range("F69").Value = "B.93"
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = "Hillerod"
$ws.Range("I69").Value = 3
$ws.Range("J69").Value = 2.6
$ws.Range("K69").Value = "30/09/2023 13:12"
$ws.Range("L69").Value = 3.6
$ws.Range("M69").Value = "06/10/2023 18:59"
$ws.Range("N69").Value = 3.5
$ws.Range("O69").Value = "30/09/2023 13:12"
$ws.Range("P69").Value = 3.81
$ws.Range("Q69").Value = "06/10/2023 18:59"
$ws.Range("R69").Value = 2.5
$ws.Range("S69").Value = "30/09/2023 13:12"
$ws.Range("T69").Value = 1.98
$ws.Range("U69").Value = "06/10/2023 18:59"
$ws.Range("V69").Value = "https://www.betexplorer.com/football/denmark/1st-division/boldklubben-1893-hillerod/hCneYENl/"

# Row 74
$ws.Range("F74").Value = "Hillerod"
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = "Fredericia"
$ws.Range("I74").Value = 1
$ws.Range("J74").Value = 2.84
$ws.Range("K74").Value = "13/10/2023 18:13"
$ws.Range("L74").Value = 3.03
$ws.Range("M74").Value = "20/10/2023 18:58"
$ws.Range("N74").Value = 3.6
$ws.Range("O74").Value = "13/10/2023 18:13"
$ws.Range("P74").Value = 3.79
$ws.Range("Q74").Value = "20/10/2023 18:58"
$ws.Range("R74").Value = 2.26
$ws.Range("S74").Value = "13/10/2023 18:13"
$ws.Range("T74").Value = 2.31
$ws.Range("U74").Value = "20/10/2023 18:35"
$ws.Range("V74").Value = "https://www.betexplorer.com/football/denmark/1st-division/hillerod-fredericia/G4z2ifUQ/"

# Row 75
$ws.Range("F75").Value = "Hobro"
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = "Helsingor"
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 1.92
$ws.Range("K75").Value = "13/10/2023 18:13"
$ws.Range("L75").Value = 1.85
$ws.Range("M75").Value = "20/10/2023 18:38"
$ws.Range("N75").Value = 3.68
$ws.Range("O75").Value = "13/10/2023 18:13"
$ws.Range("P75").Value = 3.59
$ws.Range("Q75").Value = "20/10/2023 18:39"
$ws.Range("R75").Value = 3.56
$ws.Range("S75").Value = "13/10/2023 18:13"
$ws.Range("T75").Value = 4.38
$ws.Range("U75").Value = "20/10/2023 18:36"
$ws.Range("V75").Value = "https://www.betexplorer.com/football/denmark/1st-division/hobro-helsingor/YqNTnExr/"

# Row 76
$ws.Range("F76").Value = "Sonderjyske"
$ws.Range("G76").Value = 4
$ws.Range("H76").Value = "Kolding IF"
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 1.82
$ws.Range("K76").Value = "13/10/2023 18:13"
$ws.Range("L76").Value = 1.93
$ws.Range("M76").Value = "20/10/2023 18:59"
$ws.Range("N76").Value = 4.01
$ws.Range("O76").Value = "13/10/2023 18:13"
$ws.Range("P76").Value = 3.88
$ws.Range("Q76").Value = "20/10/2023 18:59"
$ws.Range("R76").Value = 3.61
$ws.Range("S76").Value = "13/10/2023 18:13"
$ws.Range("T76").Value = 3.8
$ws.Range("U76").Value = "20/10/2023 18:39"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/denmark/1st-division/sonderjyske-kolding-if/CzMXoYil/"

# Row 80
$ws.Range("F80").Value = "Horsens"
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = "Hillerod"
$ws.Range("I80").Value = 3
$ws.Range("J80").Value = 1.84
$ws.Range("K80").Value = "20/10/2023 19:12"
$ws.Range("L80").Value = 2.17
$ws.Range("M80").Value = "27/10/2023 18:51"
$ws.Range("N80").Value = 3.9
$ws.Range("O80").Value = "20/10/2023 19:12"
$ws.Range("P80").Value = 3.59
$ws.Range("Q80").Value = "27/10/2023 18:51"
$ws.Range("R80").Value = 3.63
$ws.Range("S80").Value = "20/10/2023 19:12"
$ws.Range("T80").Value = 3.27
$ws.Range("U80").Value = "27/10/2023 18:51"
$ws.Range("V80").Value = "https://www.betexplorer.com/football/denmark/1st-division/horsens-hillerod/G4dX6ZTD/"

# Row 81
$ws.Range("F81").Value = "Hobro"
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = "Koge"
$ws.Range("I81").Value = 1
$ws.Range("J81").Value = 1.74
$ws.Range("K81").Value = "22/10/2023 16:12"
$ws.Range("L81").Value = 1.65
$ws.Range("M81").Value = "27/10/2023 18:51"
$ws.Range("N81").Value = 3.92
$ws.Range("O81").Value = "22/10/2023 16:12"
$ws.Range("P81").Value = 4.26
$ws.Range("Q81").Value = "27/10/2023 18:51"
$ws.Range("R81").Value = 4.42
$ws.Range("S81").Value = "22/10/2023 16:12"
$ws.Range("T81").Value = 4.85
$ws.Range("U81").Value = "27/10/2023 18:51"
$ws.Range("V81").Value = "https://www.betexplorer.com/football/denmark/1st-division/hobro-koge/UBGxoh7f/"

# --- Append new match rows 103 and 104 ---
$ws.Range("A102:V102").Copy() | Out-Null
$ws.Range("A103:V104").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 103
$ws.Range("A103").Value = 102
$ws.Range("B103").Value = "denmark"
$ws.Range("C103").Value = "1st-division"
$ws.Range("D103").Value = "2023-2024"
$ws.Range("E103").Value = 45261.77083333334
$ws.Range("F103").Value = "Fredericia"
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = "Vendsyssel"
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 1.93
$ws.Range("K103").Value = "26/11/2023 14:13"
$ws.Range("L103").Value = 2.17
$ws.Range("M103").Value = "01/12/2023 18:29"
$ws.Range("N103").Value = 3.92
$ws.Range("O103").Value = "26/11/2023 14:13"
$ws.Range("P103").Value = 3.55
$ws.Range("Q103").Value = "01/12/2023 18:27"
$ws.Range("R103").Value = 3.52
$ws.Range("S103").Value = "26/11/2023 14:13"
$ws.Range("T103").Value = 3.31
$ws.Range("U103").Value = "01/12/2023 18:29"
$ws.Range("V103").Value = "https://www.betexplorer.com/football/denmark/1st-division/fredericia-vendsyssel-ff/MwuEF3Co/"

# Row 104
$ws.Range("A104").Value = 103
$ws.Range("B104").Value = "denmark"
$ws.Range("C104").Value = "1st-division"
$ws.Range("D104").Value = "2023-2024"
$ws.Range("E104").Value = 45261.79166666666
$ws.Range("F104").Value = "Horsens"
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = "Hobro"
$ws.Range("I104").Value = 1
$ws.Range("J104").Value = 1.98
$ws.Range("K104").Value = "26/11/2023 15:12"
$ws.Range("L104").Value = 2.58
$ws.Range("M104").Value = "01/12/2023 18:56"
$ws.Range("N104").Value = 3.56
$ws.Range("O104").Value = "26/11/2023 15:12"
$ws.Range("P104").Value = 3.25
$ws.Range("Q104").Value = "01/12/2023 18:56"
$ws.Range("R104").Value = 3.69
$ws.Range("S104").Value = "26/11/2023 15:12"
$ws.Range("T104").Value = 2.87
$ws.Range("U104").Value = "01/12/2023 18:56"
$ws.Range("V104").Value = "https://www.betexplorer.com/football/denmark/1st-division/horsens-hobro/8fsQC1d4/"

